# Generate Report for Handoff
# For both the zh-cn and de-de localization-status sheets, rows 4-7
# (the "Ready for handoff" / "low" priority rows) are refreshed by the
# handoff-report generation: Priority flips from "low" to "ht", and the
# Latest Handoff Datetime is stamped with the new generation time.

$wb = $excel.ActiveWorkbook

$ws_zhcn = $wb.Worksheets.Item("zh-cn")
$ws_dede = $wb.Worksheets.Item("de-de")

for ($r = 4; $r -le 7; $r++) {
    $ws_zhcn.Cells.Item($r, 5).Value = "ht"
    $ws_zhcn.Cells.Item($r, 8).Value = "2016-08-16 04:31:25"

    $ws_dede.Cells.Item($r, 5).Value = "ht"
    $ws_dede.Cells.Item($r, 8).Value = "2016-08-16 04:31:32"
}
